# Generate Report for Handoff
#
# The localization-status workbook is regenerated whenever a handoff package
# is produced. This run flips the per-language "Status" cell from the old
# handback message over to "Ready for handoff", stamps the refreshed
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps, and
# re-sizes the now-narrower status columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# The new status text is noticeably shorter than the old one, so the status
# columns collapse down from their old ~30-char width to ~16.33 chars
# (Excel quantizes ColumnWidth to whole "pixel" steps, so this is the
# closest attainable width to the refreshed auto-sized column).
$statusColWidth = 16.333333333333332

# --- Overview sheet: one row summarizing each target language -----------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-30 21:07:19"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-30 21:07:14"

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-30 21:07:19"

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
